{"js": "// Apply the textual updates described by the diff:\n// 1. Update the \"Ativa\u00e7\u00e3o\" date from 2020 to 2022.\n// 2. Trim \"Ciclo Hidrol\u00f3gico; \" from the short PT program summary and drop\n//    its trailing period.\n// 3. Trim \"Hydrological Cycle; \" from the short EN program summary.\n// 4. Rewrite the long PT \"Programa\" paragraph into a dashed list.\n// 5. Rewrite the long EN \"Programa\" paragraph into a dashed list.\n\nconst replacements = [\n  {\n    find: \"Ativa\u00e7\u00e3o: 01/01/2020\",\n    replace: \"Ativa\u00e7\u00e3o: 01/01/2022\"\n  },\n  {\n    find: \"Geomorfologia Fluvial; Padr\u00f5es de Drenagem; Ciclo Hidrol\u00f3gico; Escoamentos hidr\u00e1ulicos; medidores; bocais; instrumentos de medi\u00e7\u00e3o.\",\n    replace: \"Geomorfologia Fluvial; Padr\u00f5es de Drenagem; Escoamentos hidr\u00e1ulicos; medidores; bocais; instrumentos de medi\u00e7\u00e3o\"\n  },\n  {\n    find: \"River Geomorphology; Drainage Patterns; Hydrological Cycle; Hydraulic flow; meters; nozzles; measuring instruments.\",\n    replace: \"River Geomorphology; Drainage Patterns; Hydraulic flow; meters; nozzles; measuring instruments.\"\n  },\n  {\n    find: \"As teorias geomorfol\u00f3gicas; Processos e Formas do relevo; Processos fluviais, morfologias fluviais e padr\u00f5es de drenagem; Ciclo hidrol\u00f3gico; precipita\u00e7\u00e3o; infiltra\u00e7\u00e3o; evapotranspira\u00e7\u00e3o; escoamento superficial; instrumentos de medi\u00e7\u00e3o; opera\u00e7\u00e3o de reservat\u00f3rios; vaz\u00f5es m\u00e1ximas e m\u00ednimas: distribui\u00e7\u00e3o de frequ\u00eancia, hidrograma unit\u00e1rio; Recursos H\u00eddricos e Balan\u00e7o H\u00eddrico; propaga\u00e7\u00e3o de ondas: amortecimento em reservat\u00f3rios, amortecimento em canais; Demanda de \u00e1gua e disponibilidade dos recursos h\u00eddricos. \u00c1gua subterr\u00e2nea, aqu\u00edferos e po\u00e7os; modelo matem\u00e1tico de transforma\u00e7\u00e3o de chuva-vaz\u00e3o.\",\n    replace: \"- As teorias geomorfol\u00f3gicas;- Processos e Formas do relevo;- Processos fluviais, morfologias fluviais e padr\u00f5es de drenagem;- Precipita\u00e7\u00e3o;- Infiltra\u00e7\u00e3o;- Evapotranspira\u00e7\u00e3o;- Escoamento superficial;- Instrumentos de medi\u00e7\u00e3o (Calhas, vertedores e registros);- Opera\u00e7\u00e3o de reservat\u00f3rios;- Vaz\u00f5es m\u00e1ximas e m\u00ednimas: distribui\u00e7\u00e3o de frequ\u00eancia, hidrograma unit\u00e1rio.- \u00c1gua subterr\u00e2nea, aqu\u00edferos e po\u00e7os;\"\n  },\n  {\n    find: \"Geomorphological theories; Processes and forms of relief; Fluvial processes, river morphologies and drainage patterns; Hydrological cycle; precipitation; infiltration; evapotranspiration; surface runoff; measuring instruments; operation of reservoirs; maximum and minimum flows: frequency distribution, unit hydrograph; Water Resources and Water Balance; wave propagation: damping in reservoirs, damping in channels; Water demand and availability of water resources. Groundwater, aquifers and wells; mathematical model of rain-flow transformation.\",\n    replace: \"- Geomorphological theories;- Processes and Forms of relief;- River processes, river morphologies and drainage patterns;- Precipitation;- Infiltration;- Evapotranspiration;- Surface runoff;- Measuring instruments (gutters, spillways and registers);- Reservoir operation;- Maximum and minimum flow rates: frequency distribution, unit hydrograph.- Groundwater, aquifers and wells;\"\n  }\n];\n\nfor (const { find, replace } of replacements) {\n  const results = context.document.body.search(find, { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(`Text not found: ${find}`);\n  }\n\n  for (const range of results.items) {\n    range.insertText(replace, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Apply the textual updates described by the diff:\n# 1. Update the \"Ativa\u00e7\u00e3o\" date from 2020 to 2022.\n# 2. Trim \"Ciclo Hidrol\u00f3gico; \" from the short PT program summary and drop\n#    its trailing period.\n# 3. Trim \"Hydrological Cycle; \" from the short EN program summary.\n# 4. Rewrite the long PT \"Programa\" paragraph into a dashed list.\n# 5. Rewrite the long EN \"Programa\" paragraph into a dashed list.\n\n$wdFindContinue = 1\n$wdReplaceAll = 2\n\n$d = $word.ActiveDocument\n\nfunction Replace-Text($findText, $replaceText) {\n    $range = $d.Content\n    $find = $range.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $ok = $find.Execute(\n        $findText,\n        $false,\n        $false,\n        $false,\n        $false,\n        $false,\n        $true,\n        $wdFindContinue,\n        $false,\n        $replaceText,\n        $wdReplaceAll\n    )\n    if (-not $ok) {\n        throw \"Text not found: $findText\"\n    }\n}\n\nReplace-Text \"Ativa\u00e7\u00e3o: 01/01/2020\" \"Ativa\u00e7\u00e3o: 01/01/2022\"\n\nReplace-Text \"Geomorfologia Fluvial; Padr\u00f5es de Drenagem; Ciclo Hidrol\u00f3gico; Escoamentos hidr\u00e1ulicos; medidores; bocais; instrumentos de medi\u00e7\u00e3o.\" \"Geomorfologia Fluvial; Padr\u00f5es de Drenagem; Escoamentos hidr\u00e1ulicos; medidores; bocais; instrumentos de medi\u00e7\u00e3o\"\n\nReplace-Text \"River Geomorphology; Drainage Patterns; Hydrological Cycle; Hydraulic flow; meters; nozzles; measuring instruments.\" \"River Geomorphology; Drainage Patterns; Hydraulic flow; meters; nozzles; measuring instruments.\"\n\nReplace-Text \"As teorias geomorfol\u00f3gicas; Processos e Formas do relevo; Processos fluviais, morfologias fluviais e padr\u00f5es de drenagem; Ciclo hidrol\u00f3gico; precipita\u00e7\u00e3o; infiltra\u00e7\u00e3o; evapotranspira\u00e7\u00e3o; escoamento superficial; instrumentos de medi\u00e7\u00e3o; opera\u00e7\u00e3o de reservat\u00f3rios; vaz\u00f5es m\u00e1ximas e m\u00ednimas: distribui\u00e7\u00e3o de frequ\u00eancia, hidrograma unit\u00e1rio; Recursos H\u00eddricos e Balan\u00e7o H\u00eddrico; propaga\u00e7\u00e3o de ondas: amortecimento em reservat\u00f3rios, amortecimento em canais; Demanda de \u00e1gua e disponibilidade dos recursos h\u00eddricos. \u00c1gua subterr\u00e2nea, aqu\u00edferos e po\u00e7os; modelo matem\u00e1tico de transforma\u00e7\u00e3o de chuva-vaz\u00e3o.\" \"- As teorias geomorfol\u00f3gicas;- Processos e Formas do relevo;- Processos fluviais, morfologias fluviais e padr\u00f5es de drenagem;- Precipita\u00e7\u00e3o;- Infiltra\u00e7\u00e3o;- Evapotranspira\u00e7\u00e3o;- Escoamento superficial;- Instrumentos de medi\u00e7\u00e3o (Calhas, vertedores e registros);- Opera\u00e7\u00e3o de reservat\u00f3rios;- Vaz\u00f5es m\u00e1ximas e m\u00ednimas: distribui\u00e7\u00e3o de frequ\u00eancia, hidrograma unit\u00e1rio.- \u00c1gua subterr\u00e2nea, aqu\u00edferos e po\u00e7os;\"\n\nReplace-Text \"Geomorphological theories; Processes and forms of relief; Fluvial processes, river morphologies and drainage patterns; Hydrological cycle; precipitation; infiltration; evapotranspiration; surface runoff; measuring instruments; operation of reservoirs; maximum and minimum flows: frequency distribution, unit hydrograph; Water Resources and Water Balance; wave propagation: damping in reservoirs, damping in channels; Water demand and availability of water resources. Groundwater, aquifers and wells; mathematical model of rain-flow transformation.\" \"- Geomorphological theories;- Processes and Forms of relief;- River processes, river morphologies and drainage patterns;- Precipitation;- Infiltration;- Evapotranspiration;- Surface runoff;- Measuring instruments (gutters, spillways and registers);- Reservoir operation;- Maximum and minimum flow rates: frequency distribution, unit hydrograph.- Groundwater, aquifers and wells;\"\n"}
